$wb = $excel.ActiveWorkbook

# Data update for 2022-12-20: increments to column I (year 2022) violent crime counts
# across the Citywide Totals, By Neighborhood summary, and individual neighborhood sheets.

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 7095
$ws.Range("I3").Value = 7307
$ws.Range("I4").Value = 1686
$ws.Range("I5").Value = 691
$ws.Range("I6").Value = 8679
$ws.Range("I7").Value = 25458

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I2").Value = 235
$ws.Range("I3").Value = 257
$ws.Range("I7").Value = 790

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I3").Value = 358
$ws.Range("I7").Value = 962

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("I2").Value = 77
$ws.Range("I7").Value = 256

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I3").Value = 178
$ws.Range("I6").Value = 181
$ws.Range("I7").Value = 593

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I2").Value = 206
$ws.Range("I6").Value = 185
$ws.Range("I7").Value = 797
$ws.Range("I11").Value = 390
$ws.Range("I19").Value = 711
$ws.Range("I20").Value = 628
$ws.Range("I21").Value = 114
$ws.Range("I27").Value = 222
$ws.Range("I28").Value = 14
$ws.Range("I29").Value = 1516
$ws.Range("I31").Value = 256
$ws.Range("I33").Value = 1121
$ws.Range("I37").Value = 790
$ws.Range("I42").Value = 957
$ws.Range("I43").Value = 220
$ws.Range("I47").Value = 183
$ws.Range("I48").Value = 324
$ws.Range("I51").Value = 295
$ws.Range("I52").Value = 576
$ws.Range("I54").Value = 497
$ws.Range("I60").Value = 146
$ws.Range("I63").Value = 82
$ws.Range("I65").Value = 593
$ws.Range("I67").Value = 962
$ws.Range("I72").Value = 100
$ws.Range("I73").Value = 230
$ws.Range("I76").Value = 365
$ws.Range("I77").Value = 157
$ws.Range("I79").Value = 732
$ws.Range("I83").Value = 554
$ws.Range("I85").Value = 1132
$ws.Range("I87").Value = 69
$ws.Range("I90").Value = 331
$ws.Range("I101").Value = 25458

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I2").Value = 184
$ws.Range("I6").Value = 123
$ws.Range("I7").Value = 554

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I6").Value = 361
$ws.Range("I7").Value = 1121

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I6").Value = 241
$ws.Range("I7").Value = 497

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 447
$ws.Range("I3").Value = 515
$ws.Range("I4").Value = 82
$ws.Range("I6").Value = 423
$ws.Range("I7").Value = 1516

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I2").Value = 230
$ws.Range("I6").Value = 227
$ws.Range("I7").Value = 711

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("I6").Value = 165
$ws.Range("I7").Value = 324

$ws = $wb.Worksheets.Item("River North")
$ws.Range("I6").Value = 169
$ws.Range("I7").Value = 365

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I4").Value = 52
$ws.Range("I7").Value = 1132

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("I2").Value = 68
$ws.Range("I7").Value = 185

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I3").Value = 274
$ws.Range("I7").Value = 957

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("I6").Value = 86
$ws.Range("I7").Value = 114

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I3").Value = 240
$ws.Range("I7").Value = 732

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I5").Value = 17
$ws.Range("I6").Value = 221
$ws.Range("I7").Value = 628

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I3").Value = 116
$ws.Range("I6").Value = 109

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I3").Value = 186
$ws.Range("I7").Value = 576

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("I2").Value = 44
$ws.Range("I7").Value = 183

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("I6").Value = 111
$ws.Range("I7").Value = 390

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("I3").Value = 69
$ws.Range("I7").Value = 230

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("I4").Value = 18
$ws.Range("I7").Value = 206

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("I2").Value = 62
$ws.Range("I7").Value = 222

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I2").Value = 105
$ws.Range("I7").Value = 331

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I3").Value = 80
$ws.Range("I7").Value = 295

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("I2").Value = 53
$ws.Range("I7").Value = 146

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("I2").Value = 49
$ws.Range("I7").Value = 220

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("I3").Value = 21
$ws.Range("I7").Value = 100

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("I3").Value = 53
$ws.Range("I6").Value = 38
$ws.Range("I7").Value = 157

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I2").Value = 262
$ws.Range("I7").Value = 797

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("I2").Value = 10
$ws.Range("I7").Value = 69

$ws = $wb.Worksheets.Item("Edison Park")
$ws.Range("I3").Value = 6
$ws.Range("I7").Value = 14
